$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Delete the "id" column (A) and the "target_model" column (originally D, now C
# after the first delete) so only process / name_process_model / path_folder /
# artifact_name remain.
$ws.Range("A1").EntireColumn.Delete()
$ws.Range("C1").EntireColumn.Delete()

# Update the process / name_process_model values to the new "PR_" prefixed
# naming scheme.
$ws.Range("A2").Value = "PR_A"
$ws.Range("A3").Value = "PR_B"
$ws.Range("A4").Value = "PR_B"
$ws.Range("A5").Value = "PR_C"

$ws.Range("B2").Value = "PR_A_Y1"
$ws.Range("B3").Value = "PR_B_Y2"
$ws.Range("B4").Value = "PR_B_Y3"
$ws.Range("B5").Value = "PR_C_Y2"

# Window view / selection changes captured in the diff.
$ws.Range("A1:A1048576").Select()
